$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O4").Value = 1.67
$ws.Range("P4").Value = 2.1
$ws.Range("Q5").Value = 2.05
$ws.Range("R5").Value = 1.8
$ws.Range("M6").Value = 1.13
$ws.Range("N6").Value = 6
$ws.Range("J11").Value = 4.5
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 2.52
$ws.Range("R11").Value = 1.53
$ws.Range("S11").Value = 1.45
$ws.Range("T11").Value = 2.37
$ws.Range("W11").Value = 9.75
$ws.Range("X11").Value = 22
$ws.Range("Y11").Value = 14
$ws.Range("AA11").Value = 45
$ws.Range("AB11").Value = 55
$ws.Range("AC11").Value = 7.3
$ws.Range("AE11").Value = 17.5
$ws.Range("AG11").Value = 5.6
$ws.Range("AH11").Value = 7.8
$ws.Range("AK11").Value = 18
$ws.Range("AO11").Value = 24
$ws.Range("AP11").Value = 32
$ws.Range("AR11").Value = 175
$ws.Range("AS11").Value = 450
$ws.Range("AT11").Value = 2.35
$ws.Range("AU11").Value = 7.6
$ws.Range("AX11").Value = 9.75
$ws.Range("N13").Value = 17
$ws.Range("S14").Value = 1.57
$ws.Range("T14").Value = 2.25
$ws.Range("AC14").Value = 6.5
$ws.Range("AE14").Value = 19
$ws.Range("AP14").Value = 29
$ws.Range("AT14").Value = 2.25
$ws.Range("AU14").Value = 9
$ws.Range("G20").Value = 1.65
$ws.Range("H20").Value = 3.6
$ws.Range("I20").Value = 5.5
$ws.Range("J20").Value = 2.3
$ws.Range("K20").Value = 2.1
$ws.Range("L20").Value = 6
$ws.Range("Z20").Value = 12
$ws.Range("AA20").Value = 15
$ws.Range("AC20").Value = 8
$ws.Range("AG20").Value = 12
$ws.Range("AH20").Value = 26
$ws.Range("AI20").Value = 19
$ws.Range("AK20").Value = 51
$ws.Range("AO20").Value = 9
$ws.Range("AQ20").Value = 29
$ws.Range("AW20").Value = 7
$ws.Range("AX20").Value = 34
$ws.Range("AZ20").Value = 126
$ws.Range("G21").Value = 1.95
$ws.Range("H21").Value = 3.2
$ws.Range("I21").Value = 4.2
$ws.Range("J21").Value = 2.75
$ws.Range("L21").Value = 4.5
$ws.Range("Q21").Value = 2.35
$ws.Range("R21").Value = 1.57
$ws.Range("U21").Value = 2
$ws.Range("V21").Value = 1.73
$ws.Range("W21").Value = 6
$ws.Range("X21").Value = 8.5
$ws.Range("Y21").Value = 9
$ws.Range("Z21").Value = 17
$ws.Range("AF21").Value = 67
$ws.Range("AG21").Value = 9.5
$ws.Range("AH21").Value = 19
$ws.Range("AI21").Value = 15
$ws.Range("AK21").Value = 41
$ws.Range("AN21").Value = 3.75
$ws.Range("AO21").Value = 11
$ws.Range("AP21").Value = 23
$ws.Range("AU21").Value = 9
$ws.Range("AX21").Value = 23
$ws.Range("AZ21").Value = 81
$ws.Range("BA21").Value = 126
$ws.Range("BB21").Value = 301
$ws.Range("M25").Value = 1.06
$ws.Range("N25").Value = 10
$ws.Range("O25").Value = 1.3
$ws.Range("P25").Value = 3.4
$ws.Range("Q25").Value = 2.03
$ws.Range("R25").Value = 1.83
$ws.Range("J27").Value = 2.8
$ws.Range("L27").Value = 3.65
$ws.Range("S27").Value = 1.4
$ws.Range("T27").Value = 2.7
$ws.Range("AG27").Value = 8.75
$ws.Range("AM27").Value = 500
$ws.Range("AR27").Value = 75
$ws.Range("AT27").Value = 2.7
$ws.Range("AU27").Value = 7
$ws.Range("J28").Value = 3
$ws.Range("L28").Value = 3.3
$ws.Range("O28").Value = 1.32
$ws.Range("P28").Value = 3.1
$ws.Range("W28").Value = 8.5
$ws.Range("X28").Value = 13
$ws.Range("Y28").Value = 9.25
$ws.Range("Z28").Value = 28
$ws.Range("AA28").Value = 20
$ws.Range("AB28").Value = 28
$ws.Range("AD28").Value = 6.1
$ws.Range("AE28").Value = 13
$ws.Range("AG28").Value = 8.5
$ws.Range("AH28").Value = 13.5
$ws.Range("AI28").Value = 10
$ws.Range("AK28").Value = 24
$ws.Range("AL28").Value = 32
$ws.Range("AO28").Value = 13
$ws.Range("AP28").Value = 19
$ws.Range("AR28").Value = 80
$ws.Range("AS28").Value = 200
$ws.Range("AU28").Value = 6.7
$ws.Range("AX28").Value = 15
$ws.Range("AY28").Value = 21
$ws.Range("AZ28").Value = 65
$ws.Range("BA28").Value = 100
$ws.Range("G29").Value = 2.65
$ws.Range("H29").Value = 3.5
$ws.Range("I29").Value = 2.3
$ws.Range("J29").Value = 3.2
$ws.Range("L29").Value = 2.87
$ws.Range("W29").Value = 10
$ws.Range("X29").Value = 14.5
$ws.Range("Y29").Value = 10
$ws.Range("Z29").Value = 30
$ws.Range("AA29").Value = 21
$ws.Range("AB29").Value = 27
$ws.Range("AD29").Value = 6.9
$ws.Range("AF29").Value = 50
$ws.Range("AG29").Value = 9.25
$ws.Range("AH29").Value = 12
$ws.Range("AI29").Value = 9.25
$ws.Range("AJ29").Value = 23
$ws.Range("AK29").Value = 18
$ws.Range("AL29").Value = 25
$ws.Range("AN29").Value = 4.75
$ws.Range("AO29").Value = 14
$ws.Range("AP29").Value = 20
$ws.Range("AQ29").Value = 60
$ws.Range("AR29").Value = 90
$ws.Range("AU29").Value = 6.9
$ws.Range("AW29").Value = 4.4
$ws.Range("AX29").Value = 12
$ws.Range("AY29").Value = 19
$ws.Range("AZ29").Value = 45
$ws.Range("BA29").Value = 75
$ws.Range("BB29").Value = 200
$ws.Range("G30").Value = 1.88
$ws.Range("I30").Value = 3.3
$ws.Range("J30").Value = 2.37
$ws.Range("L30").Value = 3.6
$ws.Range("P30").Value = 4.7
$ws.Range("S30").Value = 1.26
$ws.Range("T30").Value = 3.5
$ws.Range("W30").Value = 11.25
$ws.Range("X30").Value = 12
$ws.Range("Z30").Value = 18
$ws.Range("AB30").Value = 18
$ws.Range("AG30").Value = 16
$ws.Range("AH30").Value = 23
$ws.Range("AJ30").Value = 45
$ws.Range("AK30").Value = 25
$ws.Range("AL30").Value = 24
$ws.Range("AO30").Value = 9.25
$ws.Range("AP30").Value = 14
$ws.Range("AQ30").Value = 29
$ws.Range("AR30").Value = 45
$ws.Range("AT30").Value = 3.5
$ws.Range("AW30").Value = 5.8
$ws.Range("AX30").Value = 17
$ws.Range("AY30").Value = 18.5
$ws.Range("AZ30").Value = 70
$ws.Range("G33").Value = 2.87
$ws.Range("I33").Value = 2.12
$ws.Range("J33").Value = 3.4
$ws.Range("L33").Value = 2.7
$ws.Range("W33").Value = 10.5
$ws.Range("X33").Value = 16
$ws.Range("Y33").Value = 10.5
$ws.Range("AA33").Value = 23
$ws.Range("AB33").Value = 29
$ws.Range("AH33").Value = 11
$ws.Range("AI33").Value = 9
$ws.Range("AJ33").Value = 20
$ws.Range("AL33").Value = 25
$ws.Range("AN33").Value = 5
$ws.Range("AO33").Value = 15
$ws.Range("AP33").Value = 21
$ws.Range("AQ33").Value = 65
$ws.Range("AR33").Value = 90
$ws.Range("AU33").Value = 6.9
$ws.Range("AW33").Value = 4.2
$ws.Range("AX33").Value = 10.75
$ws.Range("I35").Value = 4.3
$ws.Range("J35").Value = 2.45
$ws.Range("K35").Value = 2.02
$ws.Range("Q35").Value = 1.91
$ws.Range("T35").Value = 2.47
$ws.Range("V35").Value = 1.91
$ws.Range("W35").Value = 7
$ws.Range("X35").Value = 9
$ws.Range("AC35").Value = 9
$ws.Range("AG35").Value = 12
$ws.Range("AH35").Value = 25
$ws.Range("AJ35").Value = 75
$ws.Range("AL35").Value = 45
$ws.Range("AM35").Value = 500
$ws.Range("AN35").Value = 3.7
$ws.Range("AP35").Value = 18
$ws.Range("AR35").Value = 65
$ws.Range("AT35").Value = 2.45
$ws.Range("AV35").Value = 60
$ws.Range("BB35").Value = 400
